$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '274.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.55%'
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '13'
$ws.Range("G2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.63%'
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '13'
$ws.Range("G3").Style = "Normal"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.860'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.34%'
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '13'
$ws.Range("G4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06297'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.74%'
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '13'
$ws.Range("G5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.929'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.18%'
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '13'
$ws.Range("G6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.341'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '39.96%'
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '13'
$ws.Range("G7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8686'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.28%'
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '13'
$ws.Range("G8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1555'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '4.88%'
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '13'
$ws.Range("G9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05019'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.89%'
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '13'
$ws.Range("G10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07440'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.62%'
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '13'
$ws.Range("G11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02952'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.76%'
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '13'
$ws.Range("G12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09042'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.25%'
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '13'
$ws.Range("G13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001575'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.01%'
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '13'
$ws.Range("G14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006354'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.20%'
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '13'
$ws.Range("G15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005993'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.71%'
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '13'
$ws.Range("G16").Style = "Normal"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.14%'
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '13'
$ws.Range("G17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.308'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-2.13%'
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '13'
$ws.Range("G18").Style = "Normal"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.51%'
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '13'
$ws.Range("G19").Style = "Normal"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.87%'
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '13'
$ws.Range("G20").Style = "Normal"

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.06%'
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '13'
$ws.Range("G21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.912'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.49%'
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '13'
$ws.Range("G22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04382'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.50%'
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '13'
$ws.Range("G23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001169'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.00%'
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '13'
$ws.Range("G24").Style = "Normal"

# Row 25
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '13'
$ws.Range("G25").Style = "Normal"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001201'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.01%'
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '13'
$ws.Range("G26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001617'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-4.33%'
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '13'
$ws.Range("G27").Style = "Normal"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '13'
$ws.Range("G28").Style = "Normal"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '13'
$ws.Range("G29").Style = "Normal"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '13'
$ws.Range("G30").Style = "Normal"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '13'
$ws.Range("G31").Style = "Normal"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '13'
$ws.Range("G32").Style = "Normal"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '13'
$ws.Range("G33").Style = "Normal"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '13'
$ws.Range("G34").Style = "Normal"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '13'
$ws.Range("G35").Style = "Normal"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '13'
$ws.Range("G36").Style = "Normal"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '13'
$ws.Range("G37").Style = "Normal"

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '13'
$ws.Range("G38").Style = "Normal"

# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '13'
$ws.Range("G39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04060'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.67%'
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '13'
$ws.Range("G40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007050'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4.74%'
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '13'
$ws.Range("G41").Style = "Normal"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1170'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.61%'
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '13'
$ws.Range("G42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002071'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.90%'
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '13'
$ws.Range("G43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01118'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-10.67%'
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '13'
$ws.Range("G44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005218'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.08%'
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '13'
$ws.Range("G45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = 'BOLO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.490'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-37.34%'
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '13'
$ws.Range("G46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.02000'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-11.20%'
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '13'
$ws.Range("G47").Style = "Normal"

# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '13'
$ws.Range("G48").Style = "Normal"

# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '13'
$ws.Range("G49").Style = "Normal"

# Row 50
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '13'
$ws.Range("G50").Style = "Normal"

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '13'
$ws.Range("G51").Style = "Normal"
